# Add the missing "2.11.2018" work-log entry to the "Eetu Pihamäki" time
# tracking sheet (row 25 of the Table3 listing): date, start/end time and
# the task description. The "Työaika" (duration) column is a shared
# formula already present on the row, so it recalculates on its own, as
# do the SUM()-based totals on this sheet and on "Summasivu".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Pvm (date) - 2.11.2018
$ws.Range("A25").Value = 43406
# Aloitusklo (start time) - 9:30
$ws.Range("B25").Value = 0.39583333333333331
# Lopetusklo (end time) - 16:15
$ws.Range("C25").Value = 0.67708333333333337
# Tehtävä (task description)
$ws.Range("F25").Value = "1h etsin lokit midPointista ja exportasin ne csv-tiedostona. Lisäsin ne GitHubiin kansion ""midPoint lokit"" alle. 5 h asensin mm. eclipsen eri versioita Ubuntu Desktop VM:ään ja yritin saada toimimaan. Jäi kesken. Piti saada Log Viewer toimimaan. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%202.11.2018.txt"

# The row grew taller to fit the wrapped task description.
$ws.Rows.Item(25).RowHeight = 105

# Leave the selection where the author ended up after typing the entry.
$ws.Activate() | Out-Null
$ws.Range("F25").Select() | Out-Null
